# Implemented getting kafka relations.
# The classFields sheet rows were regenerated in a different (non-alphabetical,
# re-scanned) order. Re-write the classFields rows (2-19) in their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$values = @(
    @("org.andante.config.security.role.KeycloakRole", "BLOGGER", "public", "org.andante.config.security.role.KeycloakRole"),
    @("org.andante.config.security.role.KeycloakRole", "name", "private", "java.lang.String"),
    @("org.andante.config.security.role.KeycloakRole", "`$VALUES", "private", "org.andante.config.security.role.KeycloakRole[]"),
    @("org.andante.config.security.role.KeycloakRole", "ADMIN", "public", "org.andante.config.security.role.KeycloakRole"),
    @("org.andante.config.gateway.GatewayConfiguration", "filterFactory", "private", "org.springframework.cloud.gateway.filter.factory.TokenRelayGatewayFilterFactory"),
    @("org.andante.config.security.filter.CrossOriginRequestSharingFilter", "allowedHeaders", "private", "java.lang.String"),
    @("org.andante.config.security.filter.CrossOriginRequestSharingFilter", "exposedHeaders", "private", "java.lang.String"),
    @("org.andante.config.security.filter.CrossOriginRequestSharingFilter", "allowedOrigins", "private", "java.lang.String"),
    @("org.andante.config.security.filter.CrossOriginRequestSharingFilter", "allowedMethods", "private", "java.lang.String"),
    @("org.andante.config.security.converter.KeycloakRealmRoleConverter", "ROLES", "private", "java.lang.String"),
    @("org.andante.config.security.converter.KeycloakRealmRoleConverter", "REALM_ACCESS", "private", "java.lang.String"),
    @("org.andante.config.security.SecurityConfiguration", "keycloakRealmRoleConverter", "private", "org.andante.config.security.converter.KeycloakRealmRoleConverter"),
    @("org.andante.config.security.SecurityConfiguration", "allowedHeaders", "private", "java.util.List"),
    @("org.andante.config.security.SecurityConfiguration", "jwkSetUri", "private", "java.lang.String"),
    @("org.andante.config.security.SecurityConfiguration", "allowedOrigins", "private", "java.util.List"),
    @("org.andante.config.security.SecurityConfiguration", "exposedHeaders", "private", "java.util.List"),
    @("org.andante.config.security.SecurityConfiguration", "disabledSecurityEndpoints", "private", "java.util.List"),
    @("org.andante.config.security.SecurityConfiguration", "allowedMethods", "private", "java.util.List")
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $rowData = $values[$i]
    $ws.Range("A$row").Value = $rowData[0]
    $ws.Range("B$row").Value = $rowData[1]
    $ws.Range("C$row").Value = $rowData[2]
    $ws.Range("D$row").Value = $rowData[3]
}
